$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '25.915.43'

$ws.Cells.Item(2, 5).Value = '  -0.62%  '

$ws.Cells.Item(3, 4).Value = '1.743.23'

$ws.Cells.Item(3, 5).Value = '  -1.23%  '

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '0.9995'
$ws.Cells.Item(4, 4).Style = "Normal"

$ws.Cells.Item(4, 5).Value = '  -0.22%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '231.58'
$ws.Cells.Item(5, 4).Style = "Normal"

$ws.Cells.Item(5, 5).Value = '  -2.61%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '0.9993'
$ws.Cells.Item(6, 4).Style = "Normal"

$ws.Cells.Item(6, 5).Value = '  -0.19%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.5264'
$ws.Cells.Item(7, 4).Style = "Normal"

$ws.Cells.Item(7, 5).Value = '  +0.49%  '

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.2771'
$ws.Cells.Item(8, 4).Style = "Normal"

$ws.Cells.Item(8, 5).Value = '  +0.48%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '39.59'
$ws.Cells.Item(9, 4).Style = "Normal"

$ws.Cells.Item(9, 5).Value = '  -2.03%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.06136'
$ws.Cells.Item(10, 4).Style = "Normal"

$ws.Cells.Item(10, 5).Value = '  -1.13%  '

$ws.Cells.Item(11, 4).Value = '1.737.06'

$ws.Cells.Item(11, 5).Value = '  -2.29%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.07114'
$ws.Cells.Item(12, 4).Style = "Normal"

$ws.Cells.Item(12, 5).Value = '  +1.43%  '

$ws.Cells.Item(13, 5).Value = '  -3.13%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '0.6465'
$ws.Cells.Item(14, 4).Style = "Normal"

$ws.Cells.Item(14, 5).Value = '  +0.77%  '

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '4.532'
$ws.Cells.Item(15, 4).Style = "Normal"

$ws.Cells.Item(15, 5).Value = '  -0.19%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '77.31'
$ws.Cells.Item(16, 4).Style = "Normal"

$ws.Cells.Item(16, 5).Value = '  -1.24%  '

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '0.9993'
$ws.Cells.Item(17, 4).Style = "Normal"

$ws.Cells.Item(17, 5).Value = '  -0.19%  '

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.9992'
$ws.Cells.Item(18, 4).Style = "Normal"

$ws.Cells.Item(18, 5).Value = '  -0.18%  '

$ws.Cells.Item(19, 4).Value = '25.882.54'

$ws.Cells.Item(19, 5).Value = '  -0.81%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '11.57'
$ws.Cells.Item(20, 4).Style = "Normal"

$ws.Cells.Item(20, 5).Value = '  -0.81%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '0.000006686'
$ws.Cells.Item(21, 4).Style = "Normal"

$ws.Cells.Item(21, 5).Value = '  -1.12%  '

$ws.Cells.Item(22, 4).Value = '1.958.80'

$ws.Cells.Item(22, 5).Value = '  -2.15%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '4.272'
$ws.Cells.Item(23, 4).Style = "Normal"

$ws.Cells.Item(23, 5).Value = '  +4.74%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '8.786'
$ws.Cells.Item(24, 4).Style = "Normal"

$ws.Cells.Item(24, 5).Value = '  +3.92%  '

$ws.Cells.Item(25, 5).Value = '  -0.34%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '140.66'
$ws.Cells.Item(26, 4).Style = "Normal"

$ws.Cells.Item(26, 5).Value = '  +1.15%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '1.521'
$ws.Cells.Item(27, 4).Style = "Normal"

$ws.Cells.Item(27, 5).Value = '  +0.74%  '

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '15.25'
$ws.Cells.Item(28, 4).Style = "Normal"

$ws.Cells.Item(28, 5).Value = '  +0.37%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '1.809'
$ws.Cells.Item(29, 4).Style = "Normal"

$ws.Cells.Item(29, 5).Value = '  -2.05%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '102.86'
$ws.Cells.Item(30, 4).Style = "Normal"

$ws.Cells.Item(30, 5).Value = '  -0.50%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '0.08360'
$ws.Cells.Item(31, 4).Style = "Normal"

$ws.Cells.Item(31, 5).Value = '  -0.61%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '3.748'
$ws.Cells.Item(32, 4).Style = "Normal"

$ws.Cells.Item(32, 5).Value = '  +1.16%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '3.574'
$ws.Cells.Item(33, 4).Style = "Normal"

$ws.Cells.Item(33, 5).Value = '  +3.32%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.04530'
$ws.Cells.Item(34, 4).Style = "Normal"

$ws.Cells.Item(34, 5).Value = '  +1.59%  '

$ws.Cells.Item(35, 5).Value = '  -0.27%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.9800'
$ws.Cells.Item(36, 4).Style = "Normal"

$ws.Cells.Item(36, 5).Value = '  -2.55%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.6239'
$ws.Cells.Item(37, 4).Style = "Normal"

$ws.Cells.Item(37, 5).Value = '  +2.69%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '2.701'
$ws.Cells.Item(38, 4).Style = "Normal"

$ws.Cells.Item(38, 5).Value = '  -1.65%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.01591'
$ws.Cells.Item(39, 4).Style = "Normal"

$ws.Cells.Item(39, 5).Value = '  -0.05%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '1.926'
$ws.Cells.Item(40, 4).Style = "Normal"

$ws.Cells.Item(40, 5).Value = '  -3.21%  '

$ws.Cells.Item(41, 5).Value = '  -0.29%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '100.44'
$ws.Cells.Item(42, 4).Style = "Normal"

$ws.Cells.Item(42, 5).Value = '  -2.24%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.3877'
$ws.Cells.Item(43, 4).Style = "Normal"

$ws.Cells.Item(43, 5).Value = '  -0.25%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.7371'
$ws.Cells.Item(44, 4).Style = "Normal"

$ws.Cells.Item(44, 5).Value = '  -1.05%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '5.036'
$ws.Cells.Item(45, 4).Style = "Normal"

$ws.Cells.Item(45, 5).Value = '  +1.89%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.05335'
$ws.Cells.Item(46, 4).Style = "Normal"

$ws.Cells.Item(46, 5).Value = '  -3.28%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.1127'
$ws.Cells.Item(47, 4).Style = "Normal"

$ws.Cells.Item(47, 5).Value = '  +0.53%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '6.251'
$ws.Cells.Item(48, 4).Style = "Normal"

$ws.Cells.Item(48, 5).Value = '  -1.58%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '53.84'
$ws.Cells.Item(49, 4).Style = "Normal"

$ws.Cells.Item(49, 5).Value = '  +2.18%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '30.20'
$ws.Cells.Item(50, 4).Style = "Normal"

$ws.Cells.Item(50, 5).Value = '  -0.09%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '7.685'
$ws.Cells.Item(51, 4).Style = "Normal"

$ws.Cells.Item(51, 5).Value = '  +2.73%  '
